$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '52.353.89'
$ws.Range('E2').Value = '  +5.81%  '
$ws.Range('D3').Value = '2.799.11'
$ws.Range('E3').Value = '  +6.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '117.33'
$ws.Range('E5').Value = '  +4.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '341.37'
$ws.Range('E6').Value = '  +4.73%  '
$ws.Range('E7').Value = '  +5.50%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +6.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.23'
$ws.Range('E10').Value = '  +6.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0872'
$ws.Range('E11').Value = '  +7.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.12'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.66'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '3.234.99'
$ws.Range('E15').Value = '  +6.00%  '
$ws.Range('D16').Value = '2.782.90'
$ws.Range('E16').Value = '  +5.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.891'
$ws.Range('E17').Value = '  +4.40%  '
$ws.Range('D18').Value = '52.152.20'
$ws.Range('E18').Value = '  +5.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.24'
$ws.Range('E19').Value = '  +11.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.43'
$ws.Range('E20').Value = '  +2.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.96'
$ws.Range('E21').Value = '  +4.37%  '
$ws.Range('D22').Value = '0.0₃0991'
$ws.Range('E22').Value = '  +4.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '279.49'
$ws.Range('E23').Value = '  +4.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.41'
$ws.Range('E24').Value = '  +1.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.82'
$ws.Range('E25').Value = '  +10.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.01'
$ws.Range('E26').Value = '  +3.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.24'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('E30').Value = '  +2.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.03'
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.44'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.76'
$ws.Range('E33').Value = '  +4.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0828'
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.13'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.10'
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.01'
$ws.Range('E38').Value = '  +1.01%  '
$ws.Range('E39').Value = '  +5.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.75'
$ws.Range('E40').Value = '  +28.95%  '
$ws.Range('E41').Value = '  +12.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.56'
$ws.Range('E42').Value = '  +3.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.34'
$ws.Range('E43').Value = '  +4.94%  '
$ws.Range('E44').Value = '  +4.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '126.43'
$ws.Range('E45').Value = '  -1.84%  '
$ws.Range('D46').Value = '2.110.74'
$ws.Range('E46').Value = '  +2.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.36'
$ws.Range('E47').Value = '  +3.21%  '
$ws.Range('E48').Value = '  +3.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.58'
$ws.Range('E49').Value = '  +7.13%  '
$ws.Range('E50').Value = '  +23.07%  '
$ws.Range('E51').Value = '  +1.20%  '
